$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns keep their exact string representation
# (avoids Excel auto-converting numeric-looking strings like "231.70" into numbers)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.455.31'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +4.96%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.249.25'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +4.00%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.70'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.56%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.638'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.33%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '64.03'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.410'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +3.32%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '59.38'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.27%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0902'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +5.14%  '

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.93%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.584.44'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +4.02%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.19'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.02%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.61'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.99%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.829'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.86%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.66'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.21%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.254.88'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +5.03%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '41.373.41'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.87%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '73.69'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +2.43%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0920'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +7.80%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.69%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '250.85'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +9.12%  '

$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.33%  '

$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.34'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.08%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.87'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.07%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '173.15'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.41%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +3.28%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.45'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.74%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.73%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.82'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +9.17%  '

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +2.27%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.06'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +7.13%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +2.98%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0635'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.07%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.99'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.42%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.85'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +7.99%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.88%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.000268'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +73.26%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.35%  '

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +14.19%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0240'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.48%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.84'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +14.27%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.56'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.76'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.42%  '

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.68%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.512.27'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.18%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0945'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.32%  '

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.88%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.79'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.28%  '
